$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 536.51514
$ws.Range("J17").Value = 547.03125
$ws.Range("L17").Value = 1641.09375
$ws.Range("N17").Value = -1977.09375

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2424.4092
$ws.Range("I138").Value = 1433.7273
$ws.Range("J138").Value = 3415.0908
$ws.Range("K138").Value = 4301.1819
$ws.Range("L138").Value = 10245.2724
$ws.Range("M138").Value = 838.8181000000004
$ws.Range("N138").Value = -20525.2724

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2126.3704
$ws.Range("I2").Value = 1192
$ws.Range("J2").Value = 4345.5
$ws.Range("K2").Value = 1192
$ws.Range("L2").Value = 4345.5
$ws.Range("M2").Value = -1079
$ws.Range("N2").Value = -4571.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 8817.75
$ws.Range("I31").Value = 6757
$ws.Range("J31").Value = 15000
$ws.Range("K31").Value = 6757
$ws.Range("L31").Value = 15000
$ws.Range("M31").Value = -6463
$ws.Range("N31").Value = -15588

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1810.2142
$ws.Range("I61").Value = 1586
$ws.Range("J61").Value = 2213.8
$ws.Range("K61").Value = 1586
$ws.Range("L61").Value = 2213.8
$ws.Range("M61").Value = -1374
$ws.Range("N61").Value = -2637.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 937.7143
$ws.Range("I74").Value = 1004.46155
$ws.Range("J74").Value = 744.8889
$ws.Range("K74").Value = 1004.46155
$ws.Range("L74").Value = 744.8889
$ws.Range("M74").Value = -130.46155
$ws.Range("N74").Value = -2492.8889

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 937.7143
$ws.Range("I77").Value = 1004.46155
$ws.Range("J77").Value = 744.8889
$ws.Range("K77").Value = 5022.30775
$ws.Range("L77").Value = 3724.4445
$ws.Range("M77").Value = -654.3077499999999
$ws.Range("N77").Value = -12460.4445

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 2705.05
$ws.Range("I110").Value = 2894.5
$ws.Range("J110").Value = 1000
$ws.Range("K110").Value = 2894.5
$ws.Range("L110").Value = 1000
$ws.Range("M110").Value = -849.5
$ws.Range("N110").Value = -5090

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 2126.3704
$ws.Range("I116").Value = 1192
$ws.Range("J116").Value = 4345.5
$ws.Range("K116").Value = 1192
$ws.Range("L116").Value = 4345.5
$ws.Range("M116").Value = 1102
$ws.Range("N116").Value = -8933.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 5129.4546
$ws.Range("I132").Value = 7442.4
$ws.Range("J132").Value = 3202
$ws.Range("K132").Value = 22327.2
$ws.Range("L132").Value = 9606
$ws.Range("M132").Value = -19797.2
$ws.Range("N132").Value = -14666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1810.2142
$ws.Range("I136").Value = 1586
$ws.Range("J136").Value = 2213.8
$ws.Range("K136").Value = 4758
$ws.Range("L136").Value = 6641.400000000001
$ws.Range("M136").Value = -2208
$ws.Range("N136").Value = -11741.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2126.3704
$ws.Range("I3").Value = 1192
$ws.Range("J3").Value = 4345.5
$ws.Range("K3").Value = 1192
$ws.Range("L3").Value = 4345.5
$ws.Range("M3").Value = -1078
$ws.Range("N3").Value = -4573.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H102").Value = 25000
$ws.Range("J102").Value = 25000
$ws.Range("L102").Value = 25000
$ws.Range("N102").Value = -31490

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H11").Value = 16500
$ws.Range("J11").Value = 16500
$ws.Range("L11").Value = 16500
$ws.Range("N11").Value = -16780

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 2108.4
$ws.Range("I114").Value = 258.30768
$ws.Range("J114").Value = 4112.6665
$ws.Range("K114").Value = 774.92304
$ws.Range("L114").Value = 12337.9995
$ws.Range("M114").Value = 2479.07696
$ws.Range("N114").Value = -18845.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 3166.7693
$ws.Range("J117").Value = 3739.3333
$ws.Range("L117").Value = 11217.9999
$ws.Range("N117").Value = -18101.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 801.93445
$ws.Range("I131").Value = 321.42856
$ws.Range("J131").Value = 945.06384
$ws.Range("K131").Value = 964.28568
$ws.Range("L131").Value = 2835.19152
$ws.Range("M131").Value = 4075.71432
$ws.Range("N131").Value = -12915.19152

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3105.611
$ws.Range("I80").Value = 2904
$ws.Range("J80").Value = 3357.625
$ws.Range("K80").Value = 2904
$ws.Range("L80").Value = 3357.625
$ws.Range("M80").Value = -1906
$ws.Range("N80").Value = -5353.625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 3105.611
$ws.Range("I83").Value = 2904
$ws.Range("J83").Value = 3357.625
$ws.Range("K83").Value = 14520
$ws.Range("L83").Value = 16788.125
$ws.Range("M83").Value = -9528
$ws.Range("N83").Value = -26772.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 695.28
$ws.Range("I107").Value = 471.55554
$ws.Range("K107").Value = 471.55554
$ws.Range("M107").Value = 1448.44446

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1711.2916
$ws.Range("I113").Value = 1678.9524
$ws.Range("J113").Value = 1937.6666
$ws.Range("K113").Value = 1678.9524
$ws.Range("L113").Value = 1937.6666
$ws.Range("M113").Value = 491.0476000000001
$ws.Range("N113").Value = -6277.6666

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2842.913
$ws.Range("I132").Value = 2255
$ws.Range("J132").Value = 3484.2727
$ws.Range("K132").Value = 6765
$ws.Range("L132").Value = 10452.8181
$ws.Range("M132").Value = -4235
$ws.Range("N132").Value = -15512.8181

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2119.1667
$ws.Range("I61").Value = 1728.5714
$ws.Range("J61").Value = 3486.25
$ws.Range("K61").Value = 1728.5714
$ws.Range("L61").Value = 3486.25
$ws.Range("M61").Value = -1526.5714
$ws.Range("N61").Value = -3890.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 2119.1667
$ws.Range("I113").Value = 1728.5714
$ws.Range("J113").Value = 3486.25
$ws.Range("K113").Value = 1728.5714
$ws.Range("L113").Value = 3486.25
$ws.Range("M113").Value = 441.4286
$ws.Range("N113").Value = -7826.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2157.6206
$ws.Range("I132").Value = 1607.625
$ws.Range("K132").Value = 4822.875
$ws.Range("M132").Value = -2292.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H133").Value = 52001
$ws.Range("J133").Value = 52001
$ws.Range("L133").Value = 52001
$ws.Range("N133").Value = -57061

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1938.25
$ws.Range("I136").Value = 1486.0714
$ws.Range("J136").Value = 2993.3333
$ws.Range("K136").Value = 4458.2142
$ws.Range("L136").Value = 8979.999899999999
$ws.Range("M136").Value = -1908.2142
$ws.Range("N136").Value = -14079.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H140").Value = 56477.332
$ws.Range("J140").Value = 56477.332
$ws.Range("L140").Value = 56477.332
$ws.Range("N140").Value = -66837.33199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H141").Value = 69750
$ws.Range("J141").Value = 69750
$ws.Range("L141").Value = 69750
$ws.Range("N141").Value = -80110

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1802.9814
$ws.Range("I132").Value = 1551.0731
$ws.Range("J132").Value = 2597.4614
$ws.Range("K132").Value = 4653.219300000001
$ws.Range("L132").Value = 7792.3842
$ws.Range("M132").Value = -2123.219300000001
$ws.Range("N132").Value = -12852.3842

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1909.72
$ws.Range("I136").Value = 1638.15
$ws.Range("J136").Value = 2996
$ws.Range("K136").Value = 4914.450000000001
$ws.Range("L136").Value = 8988
$ws.Range("M136").Value = -2364.450000000001
$ws.Range("N136").Value = -14088
